# Natmi following Dr Hou advice
# Update the Efna2-Epha2 LR-pairs sheet: existing rows 2-4 get refreshed
# values and five new sending/target cluster combinations (rows 5-10) are
# appended, covering the full ECs / FAPs / sCs cross-product.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D are cluster/gene labels, E-T are the numeric NATMI stats.
$data = @(
    @("ECs", "Efna2", "Epha2", "ECs", 1, 0.3333333333333333, 1.192082333333333, 3.576247, 0.2797939869571494, 0.2797939869571493, 3, 1, 11.08503166666667, 33.255095, 0.4259149152633459, 0.4259149152633459, 13.21427041427389, 118.928433728465, 0.119168432246048, 0.1191684322460479),
    @("ECs", "Efna2", "Epha2", "FAPs", 1, 0.3333333333333333, 1.192082333333333, 3.576247, 0.2797939869571494, 0.2797939869571493, 2, 0.6666666666666666, 0.312365, 0.937095, 0.01200185227312402, 0.01200185227312402, 0.3723647980516667, 3.351283182465, 0.003358046098368094, 0.003358046098368093),
    @("ECs", "Efna2", "Epha2", "sCs", 1, 0.3333333333333333, 1.192082333333333, 3.576247, 0.2797939869571494, 0.2797939869571493, 3, 1, 14.62900266666666, 43.88700799999999, 0.5620832324635302, 0.5620832324635302, 17.43897563321955, 156.950780698976, 0.1572675086127333, 0.1572675086127333),
    @("FAPs", "Efna2", "Epha2", "ECs", 3, 1, 1.843761666666667, 5.531285, 0.4327498305196134, 0.4327498305196134, 3, 1, 11.08503166666667, 33.255095, 0.4259149152633459, 0.4259149152633459, 20.43815646078611, 183.943408147075, 0.1843146073959884, 0.1843146073959884),
    @("FAPs", "Efna2", "Epha2", "FAPs", 3, 1, 1.843761666666667, 5.531285, 0.4327498305196134, 0.4327498305196134, 2, 0.6666666666666666, 0.312365, 0.937095, 0.01200185227312402, 0.01200185227312402, 0.5759266130083334, 5.183339517075001, 0.005193799537115854, 0.005193799537115854),
    @("FAPs", "Efna2", "Epha2", "sCs", 3, 1, 1.843761666666667, 5.531285, 0.4327498305196134, 0.4327498305196134, 3, 1, 14.62900266666666, 43.88700799999999, 0.5620832324635302, 0.5620832324635302, 26.97239433836444, 242.75154904528, 0.2432414235865092, 0.2432414235865091),
    @("sCs", "Efna2", "Epha2", "ECs", 3, 1, 1.224727666666667, 3.674183, 0.2874561825232373, 0.2874561825232373, 3, 1, 11.08503166666667, 33.255095, 0.4259149152633459, 0.4259149152633459, 13.57614496804278, 122.185304712385, 0.1224318756213095, 0.1224318756213095),
    @("sCs", "Efna2", "Epha2", "FAPs", 3, 1, 1.224727666666667, 3.674183, 0.2874561825232373, 0.2874561825232373, 2, 0.6666666666666666, 0.312365, 0.937095, 0.01200185227312402, 0.01200185227312402, 0.3825620575983333, 3.443058518385, 0.003450006637640067, 0.003450006637640067),
    @("sCs", "Efna2", "Epha2", "sCs", 3, 1, 1.224727666666667, 3.674183, 0.2874561825232373, 0.2874561825232373, 3, 1, 14.62900266666666, 43.88700799999999, 0.5620832324635302, 0.5620832324635302, 17.91654430160711, 161.248898714464, 0.1615743002642878, 0.1615743002642878)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
